$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase MaxInvest (column F) values for existing units
$ws.Range("F8").Value = 17
$ws.Range("F10").Value = 17
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 16
$ws.Range("F14").Value = 6
$ws.Range("F16").Value = 78

# Update the active selection to K22 (matches the selection change in the diff)
$ws.Range("K22").Select()
